$d = $word.ActiveDocument

# 1. "HGA" -> "HGW"
$d.Content.Find.Execute("HGA", $true, $false, $false, $false, $false, $true, 1, $false, "HGW", 2)

# 2. Split the BOSS 5 paragraph right after "...super effective " and insert the
#    new FINAL BOSS paragraph, keeping the two trailing manual line breaks that
#    originally closed the BOSS 5 paragraph.
$rng = $d.Content
$rng.Find.Execute("super effective ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertParagraphAfter()

# The paragraph that was just created (after the split) is now empty except for
# the two trailing <w:br/> runs that used to end the BOSS 5 paragraph.
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$ins = $newPara.Range
$ins.Collapse(1)
$ins.InsertBefore("FINAL BOSS: JANUSZ P(AWULON)" + [char]11 + "No i tu się zaczyna. Na wstępie musisz zniszczyć każdą z jego lasek, której bronią wszyscy Twoi poprzedni przeciwnicy. Po ich zniszczeniu wychodzi on cały na biało. Strzela do Ciebie kawałkami bliżej niezidentyfikowanego ciasta, którym się cały czas objada. Ono zadaje Ci mały dmg ale nakłada poisona. Dodatkowo co jakiś czas macha swoją boską laską. Co dokładnie 13,05 sekundy dostaje odporność na obrażenia na czas 2 sekund. Po śmierci rozbłyska jasnym światłem i znika zostawiając po sobie bulbulator, pierścienie i książkę, która zamyka się od wiatru.")

# Split "FINAL BOSS: JANUSZ P(AWULON)" and the rest into two separate runs
# (matching the target OOXML, where the line break + following text form
# their own <w:r>). A throwaway bookmark forces the run boundary without
# leaving any formatting residue behind.
$headRng = $d.Content
$headRng.Find.Execute("FINAL BOSS: JANUSZ P(AWULON)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitStart = $headRng.End
$paraEnd = $newPara.Range.End
$tailRng = $d.Range($splitStart, $paraEnd - 3)
$d.Bookmarks.Add("zzzTmpSplitMarker", $tailRng)
$d.Bookmarks.Item("zzzTmpSplitMarker").Delete()
